# Fix the wrong toolchain-menu reference in the "Add additional directories
# under ..." bullet: it pointed at
#   Project Properties -> Toolchain -> ARM/GNU Assembler -> General
# but should point at
#   Project Properties -> Toolchain -> ARM/GNU C Compiler -> Directories
# (matching the other "Project Properties -> ..." references in the
# document, which are italicised).

$d = $word.ActiveDocument

$rng = $d.Content.Duplicate
$rng.Find.ClearFormatting()
$rng.Find.Replacement.ClearFormatting()
$rng.Find.Text = "Project Properties -> Toolchain -> ARM/GNU Assembler -> General "
$rng.Find.Replacement.Text = "Project Properties -> Toolchain -> ARM/GNU C Compiler -> Directories "
$rng.Find.Replacement.Font.Italic = 1
$rng.Find.Forward = $true
$rng.Find.Wrap = 1
$rng.Find.Format = $true
$rng.Find.MatchCase = $true
$rng.Find.MatchWholeWord = $false
$rng.Find.MatchWildcards = $false

$found = $rng.Find.Execute(
    [Type]::Missing, [Type]::Missing, [Type]::Missing, [Type]::Missing,
    [Type]::Missing, [Type]::Missing, [Type]::Missing, [Type]::Missing,
    [Type]::Missing, [Type]::Missing, 2)

Write-Output "Replaced wrong toolchain reference: $found"
